# Applies the "Add files via upload" edit to Fragen-checkbox.xlsx:
#  - row 6 (col D): the "Förderungen" question string's trailing
#    ":pflicht" becomes "|pflicht"
#  - rows 7-11 (col D) on "Tabelle1": the Upload question string's trailing
#    ":pflicht" becomes "|pflicht"
#  - row 11 (col D) additionally picks up the "Text" number format (style
#    used by the sibling cells D7:D10) since it previously had no explicit
#    style
#  - row 12 (col D): the "Richtig-und-Vollständig" checkbox string's
#    trailing ":pflicht" becomes "|pflicht" (a brand-new shared string)
#  - the active selection on "Tabelle1" moves from D6 to D13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$uploadText = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"

$ws.Range("D6").Value = "Förderungen:mcheckbox(keine,AMA,ÖPUL,Sonstige);Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"

$ws.Range("D7").Value = $uploadText
$ws.Range("D8").Value = $uploadText
$ws.Range("D9").Value = $uploadText
$ws.Range("D10").Value = $uploadText

# D11 had no explicit cell style before; match D7:D10's "Text" number format.
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = $uploadText

$ws.Range("D12").Value = "Richtig-und-Vollständig:checkbox|pflicht"

# Update the selected cell on the visible sheet.
$ws.Range("D13").Select()
